$d = $word.ActiveDocument

# --- Change 1: "poskim" -> "poseqim (adjudicators of Jewish law)" ---
$d.Content.Find.Execute("poskim", $false, $false, $false, $false, $false, $true, 1, $false, "poseqim (adjudicators of Jewish law)", 2) | Out-Null

# --- Change 2: empty paragraph after "Using modern natural language..." gets a single space run ---
$pEmpty = $d.Paragraphs.Item(12)
$rEmpty = $pEmpty.Range
$rEmpty.InsertBefore(" ")
$rEmpty2 = $d.Paragraphs.Item(12).Range
$rEmpty2.Font.Name = "Calibri"
$rEmpty2.Font.NameAscii = "Calibri"
$rEmpty2.Font.NameFarEast = "Calibri"
$rEmpty2.Font.NameBi = "Calibri"
$rEmpty2.Font.NameOther = "Calibri"
$rEmpty2.Font.Size = 12
$rEmpty2.Font.SizeBi = 12

# --- Change 3: replace text of paragraphs 15,17,19,21,23,25 (keep paragraph shells + interleaved empties) ---
$p15 = $d.Paragraphs.Item(15)
$r15 = $p15.Range
$rTrim15 = $d.Range($r15.Start, $r15.End - 1)
$rTrim15.Text = "Whether doing training or inference, the model requires a tokenizer to convert the words into a vector representation of the text. This vector representation is then fed into BERT to generate an encoding. When training, that encoding is then fed to a second half of the network, which predicts the missing word, followed by a loss calculation, and the updating of weights. When doing inference, the distance between the current paragraph’s vector representations, and all other vector representations is calculated, and the closest paragraphs are represented to the user.  After being trained in this fashion, the middle layer’s output is used as a vector encoding of the input text. "

$p17 = $d.Paragraphs.Item(17)
$r17 = $p17.Range
$rTrim17 = $d.Range($r17.Start, $r17.End - 1)
$rTrim17.Text = "In April of 2021, the Bar Ilan NLP Lab released AlephBERT, a BERT model designed for the Hebrew Language. It was trained on OSCAR’s Hebrew section (a collection of articles in Hebrew from the internet), Hebrew Wikipedia, and Hebrew twitter. Because of the training data provided, the model learned very modern colloquial hebrew. However, for use in this task AlephBert would need to learn a more formal, and much older hebrew used in rabbinic responsa. "

$p19 = $d.Paragraphs.Item(19)
$r19 = $p19.Range
$rTrim19 = $d.Range($r19.Start, $r19.End - 1)
$rTrim19.Text = "This was done by fine-tuned AlephBERT via MLM on a subset of the Bar Ilan Responsa  Project’s collection of rabanic responsa. Through this process, AlephBERT became accustomed to the dialect of Hebrew used specifically in responsa."

$p21 = $d.Paragraphs.Item(21)
$r21 = $p21.Range
$rTrim21 = $d.Range($r21.Start, $r21.End - 1)
$rTrim21.Text = "The model was then able to generate embedding in vector space that accurately represented the topics of each document. Calculating the magnitude of the distance between two vectors to represent how topically different they are. For computational efficiency, this process was implemented using matrix operations, which were performed on a GPU. Now an efficient method for comparing a paragraph to all the other paragraphs in courpas of responsa and finding the paragraphs that is the most topically similar to that paragraph was possible. "

$p23 = $d.Paragraphs.Item(23)
$r23 = $p23.Range
$rTrim23 = $d.Range($r23.Start, $r23.End - 1)
$rTrim23.Text = "The results received were reasonable with this method, though the measurements are not objective.  For any given paragraph about 2 of the top 3 suggestions were deemed relevant. "

$p25 = $d.Paragraphs.Item(25)
$r25 = $p25.Range
$rTrim25 = $d.Range($r25.Start, $r25.End - 1)
$rTrim25.Text = "Generally in rabbinic responsa, the main points in the document are expressed in the beginning when the question is asked and at the end in the conclusion. So, in order to improve the search results document embeddings were also generated. This was done by using the first 250 tokens and the last 250 tokens and in the document. This size was chosen because the maximum pre-trained AlephBERT model uses is 512 tokens.  "

# --- Append two new paragraphs after paragraph 25 (now holding new text) ---
$p25b = $d.Paragraphs.Item(25)
$p25b.Range.InsertParagraphAfter()
$p26b = $d.Paragraphs.Item(26)
$p26b.Range.InsertParagraphAfter()
$p27b = $d.Paragraphs.Item(27)
$p27b.Range.InsertBefore("This method on its own generated reasonable results which were similar in quality to the previous model. After summing both scores for any given paragraph (the sum of the paragraph distance and the document distance), when then sorted by overall distance, the results were promising, though similar. There is much difficulty in assessing model improvement unless there are very large changes, noticeable in reviewing a few dozen documents. In the future, once users are available, user feedback can be a critical tool in assessing how relevant a particular article is.")

# --- Change 4: final paragraph text replacement ---
$d.Content.Find.Execute("While our search results for topical similarity were promising, there is still more work to do. Sometimes our search results contained no paragraphs that were actually similar or just one or two similar results, even though there were actually many paragraphs that were similar. This could potentially be improved by doing more epochs of fine-tuning on our model or by doing some other combination of embeddings like how we did with making the full document embeddings.  ", $false, $false, $false, $false, $false, $true, 1, $false, "Overall, while the results for topical similarity were promising, there is still more work to be done. Sometimes results contained no paragraphs that were similar or just one or two similar results, even though there were actually many paragraphs that were similar. This could potentially be improved by increasing the number of training epochs, using an ensemble of models, and/or implementing a discriminator network. ", 2) | Out-Null
